$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '63.521.20'
$ws.Range("E2").Value = '  -3.90%  '

# Row 3
$ws.Range("D3").Value = '3.123.93'
$ws.Range("E3").Value = '  -3.96%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
Set-TextValue "D5" '607.73'
$ws.Range("E5").Value = '  +0.09%  '

# Row 6
Set-TextValue "D6" '144.51'
$ws.Range("E6").Value = '  -7.89%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").Value = '3.119.71'
$ws.Range("E8").Value = '  -4.02%  '

# Row 9
Set-TextValue "D9" '0.518'
$ws.Range("E9").Value = '  -4.00%  '

# Row 10
$ws.Range("E10").Value = '  -6.58%  '

# Row 11
$ws.Range("E11").Value = '  -7.81%  '

# Row 12
Set-TextValue "D12" '0.467'
$ws.Range("E12").Value = '  -5.00%  '

# Row 13
Set-TextValue "D13" '0.0000249'
$ws.Range("E13").Value = '  -6.24%  '

# Row 14
Set-TextValue "D14" '35.16'
$ws.Range("E14").Value = '  -8.47%  '

# Row 15
$ws.Range("D15").Value = '3.631.45'
$ws.Range("E15").Value = '  -4.18%  '

# Row 16
Set-TextValue "D16" '0.116'
$ws.Range("E16").Value = '  +1.61%  '

# Row 17
$ws.Range("D17").Value = '63.603.43'
$ws.Range("E17").Value = '  -3.90%  '

# Row 18
$ws.Range("D18").Value = '3.113.33'
$ws.Range("E18").Value = '  -4.36%  '

# Row 19
Set-TextValue "D19" '6.79'
$ws.Range("E19").Value = '  -6.40%  '

# Row 20
Set-TextValue "D20" '474.66'
$ws.Range("E20").Value = '  -4.58%  '

# Row 21
Set-TextValue "D21" '14.55'
$ws.Range("E21").Value = '  -4.67%  '

# Row 22
Set-TextValue "D22" '0.702'
$ws.Range("E22").Value = '  -5.60%  '

# Row 23
Set-TextValue "D23" '7.72'
$ws.Range("E23").Value = '  -3.56%  '

# Row 24
Set-TextValue "D24" '13.48'
$ws.Range("E24").Value = '  -7.21%  '

# Row 25
Set-TextValue "D25" '83.28'
$ws.Range("E25").Value = '  -3.96%  '

# Row 26
$ws.Range("E26").Value = '  +0.29%  '

# Row 27
Set-TextValue "D27" '2.78'
$ws.Range("E27").Value = '  -7.97%  '

# Row 28
Set-TextValue "D28" '8.36'
$ws.Range("E28").Value = '  -7.41%  '

# Row 29
$ws.Range("E29").Value = '  -9.01%  '

# Row 30
$ws.Range("E30").Value = '  -3.02%  '

# Row 31
$ws.Range("E31").Value = '  -11.79%  '

# Row 32
Set-TextValue "D32" '1.00'
$ws.Range("E32").Value = '  -0.10%  '

# Row 33
Set-TextValue "D33" '2.67'
$ws.Range("E33").Value = '  -5.71%  '

# Row 34
Set-TextValue "D34" '26.15'
$ws.Range("E34").Value = '  -5.89%  '

# Row 35
Set-TextValue "D35" '1.11'
$ws.Range("E35").Value = '  -2.00%  '

# Row 36
Set-TextValue "D36" '5.92'
$ws.Range("E36").Value = '  -6.98%  '

# Row 37
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0752'
$ws.Range("E37").Value = '  -2.26%  '

# Row 38
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D38" '52.69'
$ws.Range("E38").Value = '  -5.20%  '

# Row 39
Set-TextValue "D39" '452.87'
$ws.Range("E39").Value = '  -7.86%  '

# Row 40
Set-TextValue "D40" '2.93'
$ws.Range("E40").Value = '  -13.97%  '

# Row 41
Set-TextValue "D41" '0.0391'
$ws.Range("E41").Value = '  -6.74%  '

# Row 42
Set-TextValue "D42" '0.118'
$ws.Range("E42").Value = '  -9.11%  '

# Row 43
Set-TextValue "D43" '8.30'
$ws.Range("E43").Value = '  -4.73%  '

# Row 44
$ws.Range("D44").Value = '2.836.41'
$ws.Range("E44").Value = '  -5.12%  '

# Row 45
Set-TextValue "D45" '2.27'
$ws.Range("E45").Value = '  -9.59%  '

# Row 46
Set-TextValue "D46" '0.263'
$ws.Range("E46").Value = '  -9.07%  '

# Row 47
$ws.Range("E47").Value = '  -2.09%  '

# Row 48
$ws.Range("E48").Value = '  -0.03%  '

# Row 49
Set-TextValue "D49" '26.05'
$ws.Range("E49").Value = '  -8.39%  '

# Row 50
$ws.Range("E50").Value = '  -5.18%  '

# Row 51
Set-TextValue "D51" '118.87'
$ws.Range("E51").Value = '  -2.13%  '
